$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.9999999999993117
$ws.Range("E2").Value = 0.9999999999993117

$ws.Range("D3").Value = [double]"4.531068451036246E-07"
$ws.Range("E3").Value = [double]"4.531068451036246E-07"

$ws.Range("C4").Value = $false
$ws.Range("D4").Value = 0.6173657116625142
$ws.Range("E4").Value = 0.6173657116625142

$ws.Range("D5").Value = [double]"9.273154463488359E-16"
$ws.Range("E5").Value = [double]"9.273154463488359E-16"

$ws.Range("D6").Value = [double]"5.714966185353726E-50"
$ws.Range("E6").Value = [double]"5.714966185353726E-50"

$ws.Range("F7").Value = 4.827533721923828
$ws.Range("G7").Value = 0.6666666666666666
